$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 160  # was 161; 上海·吉卜力工作室物语-沉浸式艺术展全球首站（9月-10月）
$ws.Range("F3").Value = 2412  # was 2411; 上海 洛天依歌行宇宙·无限遨游 沉浸式体验展
$ws.Range("F4").Value = 29  # was 30; 上海·【神秘的西夏陵】大空间高沉浸探险体验
$ws.Range("F6").Value = 64  # was 63; 上海·BH 零号空洞的委托   
$ws.Range("F7").Value = 277  # was 274; 上海·排球少年Only·魔都见学同人展
$ws.Range("F8").Value = 336  # was 334; 上海·第五人格同人only
$ws.Range("F9").Value = 2286  # was 2229; 上海·趣元界-INW动漫游戏展
$ws.Range("F10").Value = 1160  # was 1159; 上海·向前冲！运动番同人Only
$ws.Range("F11").Value = 1045  # was 1044; 上海·AXG彩虹领域动漫游戏嘉年华X欢迎来到彩虹岛，一个色彩斑斓的梦幻二次元世界（免费活动）
$ws.Range("F12").Value = 853  # was 851; 上海·LookLook动漫嘉年华3th
$ws.Range("F14").Value = 843  # was 842; 上海·iPR动漫-第五&原&铁&崩&零同人ONLY同好嘉年华2.0
$ws.Range("F15").Value = 1481  # was 1477; 上海·第五人格同人only-万圣狂欢宴
$ws.Range("F16").Value = 737  # was 727; 上海 ·《天官赐福》动画四周年纪念展
$ws.Range("F17").Value = 1699  # was 1696; 上海·JOing动漫游戏嘉年华
$ws.Range("F19").Value = 368  # was 360; 上海·恶作剧国乙＋代号鸢同人only
$ws.Range("F20").Value = 68  # was 67; 上海·创造力动漫游戏嘉年华2.0
$ws.Range("F21").Value = 110  # was 108; 上海·明日方舟同人ONLY
$ws.Range("F23").Value = 2626  # was 2624; 上海·iPR动漫-第五&原&铁&崩&零同人ONLY同好嘉年华

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F19").Value = 153  # was 152; 上海·majiko巡演-2024
$ws.Range("F28").Value = 183  # was 182; 上海·2024年刘明月专场生日会
$ws.Range("F29").Value = 2  # was 1; 上海·【大乐】《怦然心动·爱乐之城》奥斯卡之夜影视金曲视听音乐会
$ws.Range("F36").Value = 54  # was 52; 上海·变形金刚音乐会40周年特别版
$ws.Range("F38").Value = 346  # was 340; 上海·Ayasa LIVE TOUR 2024〜D.D.D.〜
$ws.Range("F46").Value = 298  # was 297; 上海·中村百合香粉丝见面会最终场

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 2487  # was 2484; 上海·盗墓笔记官方授权「四季同书」主题店
$ws.Range("F6").Value = 2498  # was 2496; 上海·东方明珠·「光与夜之恋 × 线条小狗 ×爱胖达文化 」线条大作战主题店
$ws.Range("F7").Value = 9538  # was 9533; 上海·大悦城·「光与夜之恋 × 线条小狗 ×爱胖达文化 」线条大作战主题餐厅
$ws.Range("F13").Value = 2787  # was 2779; 上海·2024·《世界之外》x  萌果酱谷子咖啡
$ws.Range("F14").Value = 361  # was 358; 上海·三丽鸥家族Sanrio Characters主题餐厅·海滩奇遇季
$ws.Range("F15").Value = 668  # was 666; 上海·「火影忍者疾风传 × animate cafe」

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 2487  # was 2484; 上海·盗墓笔记官方授权「四季同书」主题店
$ws.Range("F6").Value = 160  # was 161; 上海·吉卜力工作室物语-沉浸式艺术展全球首站（9月-10月）
$ws.Range("F7").Value = 2412  # was 2411; 上海 洛天依歌行宇宙·无限遨游 沉浸式体验展
$ws.Range("F8").Value = 2787  # was 2779; 上海·2024·《世界之外》x  萌果酱谷子咖啡
$ws.Range("F9").Value = 361  # was 358; 上海·三丽鸥家族Sanrio Characters主题餐厅·海滩奇遇季
$ws.Range("F11").Value = 668  # was 666; 上海·「火影忍者疾风传 × animate cafe」
$ws.Range("F16").Value = 64  # was 63; 上海·BH 零号空洞的委托   
$ws.Range("F17").Value = 277  # was 274; 上海·排球少年Only·魔都见学同人展
$ws.Range("F18").Value = 336  # was 334; 上海·第五人格同人only
$ws.Range("F20").Value = 1045  # was 1044; 上海·AXG彩虹领域动漫游戏嘉年华X欢迎来到彩虹岛，一个色彩斑斓的梦幻二次元世界（免费活动）
$ws.Range("F21").Value = 853  # was 851; 上海·LookLook动漫嘉年华3th
$ws.Range("F23").Value = 843  # was 842; 上海·iPR动漫-第五&原&铁&崩&零同人ONLY同好嘉年华2.0
$ws.Range("F28").Value = 737  # was 727; 上海 ·《天官赐福》动画四周年纪念展
$ws.Range("F31").Value = 1699  # was 1696; 上海·JOing动漫游戏嘉年华
$ws.Range("F32").Value = 368  # was 360; 上海·恶作剧国乙＋代号鸢同人only
$ws.Range("F37").Value = 183  # was 182; 上海·2024年刘明月专场生日会
$ws.Range("F38").Value = 2  # was 1; 上海·【大乐】《怦然心动·爱乐之城》奥斯卡之夜影视金曲视听音乐会
$ws.Range("F39").Value = 68  # was 67; 上海·创造力动漫游戏嘉年华2.0
$ws.Range("F41").Value = 54  # was 52; 上海·变形金刚音乐会40周年特别版
$ws.Range("F42").Value = 346  # was 340; 上海·Ayasa LIVE TOUR 2024〜D.D.D.〜
$ws.Range("F43").Value = 2626  # was 2624; 上海·iPR动漫-第五&原&铁&崩&零同人ONLY同好嘉年华
$ws.Range("F48").Value = 298  # was 297; 上海·中村百合香粉丝见面会最终场
